$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing merged cell ranges before restructuring the data
$ws.Cells.UnMerge()

# Rewrite all question/option/total rows in sorted order
$ws.Cells.Item(1, 1).Value = 'Question name'
$ws.Cells.Item(1, 2).Value = 'Option'
$ws.Cells.Item(1, 3).Value = 'Total'

$ws.Cells.Item(2, 1).Value = 'test question'
$ws.Cells.Item(2, 2).Value = 'yes'
$ws.Cells.Item(2, 3).Value = 0

$ws.Cells.Item(3, 1).ClearContents()
$ws.Cells.Item(3, 2).Value = 'no'
$ws.Cells.Item(3, 3).Value = 0

$ws.Cells.Item(4, 1).ClearContents()
$ws.Cells.Item(4, 2).Value = 'no comment'
$ws.Cells.Item(4, 3).Value = 0

$ws.Cells.Item(5, 1).Value = 'test text'
$ws.Cells.Item(5, 2).Value = 'jdkhfjdsfhj'
$ws.Cells.Item(5, 3).ClearContents()

$ws.Cells.Item(6, 1).ClearContents()
$ws.Cells.Item(6, 2).Value = 'djhfjhdsjfhdsf'
$ws.Cells.Item(6, 3).ClearContents()

$ws.Cells.Item(7, 1).ClearContents()
$ws.Cells.Item(7, 2).Value = 'dmfdsmfnsd'
$ws.Cells.Item(7, 3).ClearContents()

$ws.Cells.Item(8, 1).Value = 'New text option question'
$ws.Cells.Item(8, 2).Value = 'ffjhjhsdjhfjsdhf'
$ws.Cells.Item(8, 3).ClearContents()

$ws.Cells.Item(9, 1).ClearContents()
$ws.Cells.Item(9, 2).Value = 'djfhsdjf fdjhsdjhf'
$ws.Cells.Item(9, 3).ClearContents()

$ws.Cells.Item(10, 1).ClearContents()
$ws.Cells.Item(10, 2).Value = 'fkjkhf sdjfbsjdkhf'
$ws.Cells.Item(10, 3).ClearContents()

$ws.Cells.Item(11, 1).Value = 'hhghghghg'
$ws.Cells.Item(11, 2).Value = 'No answer'
$ws.Cells.Item(11, 3).ClearContents()

$ws.Cells.Item(12, 1).Value = 'Updated question'
$ws.Cells.Item(12, 2).Value = 'YES'
$ws.Cells.Item(12, 3).Value = 2

$ws.Cells.Item(13, 1).ClearContents()
$ws.Cells.Item(13, 2).Value = 'NO'
$ws.Cells.Item(13, 3).Value = 3

$ws.Cells.Item(14, 1).Value = 'new kfhfhkhd'
$ws.Cells.Item(14, 2).Value = 'Yes'
$ws.Cells.Item(14, 3).Value = 0

$ws.Cells.Item(15, 1).ClearContents()
$ws.Cells.Item(15, 2).Value = 'No'
$ws.Cells.Item(15, 3).Value = 0

$ws.Cells.Item(16, 1).Value = 'All multiple question 3'
$ws.Cells.Item(16, 2).Value = 'Option 1'
$ws.Cells.Item(16, 3).Value = 0

$ws.Cells.Item(17, 1).ClearContents()
$ws.Cells.Item(17, 2).Value = 'Option 2'
$ws.Cells.Item(17, 3).Value = 0

$ws.Cells.Item(18, 1).ClearContents()
$ws.Cells.Item(18, 2).Value = 'Option 3'
$ws.Cells.Item(18, 3).Value = 0

$ws.Cells.Item(19, 1).Value = 'All multiple question 2'
$ws.Cells.Item(19, 2).Value = 'Option 1'
$ws.Cells.Item(19, 3).Value = 1

$ws.Cells.Item(20, 1).ClearContents()
$ws.Cells.Item(20, 2).Value = 'Option 2'
$ws.Cells.Item(20, 3).Value = 2

$ws.Cells.Item(21, 1).ClearContents()
$ws.Cells.Item(21, 2).Value = 'Option 3'
$ws.Cells.Item(21, 3).Value = 1

$ws.Cells.Item(22, 1).Value = 'All multiple question 1'
$ws.Cells.Item(22, 2).Value = 'jdfhdjsfhjsdhf'
$ws.Cells.Item(22, 3).ClearContents()

$ws.Cells.Item(23, 1).ClearContents()
$ws.Cells.Item(23, 2).Value = 'dfjhdjsd fksjdhfjksd'
$ws.Cells.Item(23, 3).ClearContents()

$ws.Cells.Item(24, 1).ClearContents()
$ws.Cells.Item(24, 2).Value = 'mdfd fsdjbfjsdhfjhsdf'
$ws.Cells.Item(24, 3).ClearContents()

$ws.Cells.Item(25, 1).Value = 'Question two?'
$ws.Cells.Item(25, 2).Value = 'YES'
$ws.Cells.Item(25, 3).Value = 3

$ws.Cells.Item(26, 1).ClearContents()
$ws.Cells.Item(26, 2).Value = 'NO'
$ws.Cells.Item(26, 3).Value = 1

$ws.Cells.Item(27, 1).ClearContents()
$ws.Cells.Item(27, 2).Value = 'NO COMMENT'
$ws.Cells.Item(27, 3).Value = 2

$ws.Cells.Item(28, 1).ClearContents()
$ws.Cells.Item(28, 2).Value = 'Another'
$ws.Cells.Item(28, 3).Value = 0

# Apply merged cell ranges matching the new layout
$ws.Range("A2:A4").Merge()
$ws.Range("B5:C5").Merge()
$ws.Range("B6:C6").Merge()
$ws.Range("B7:C7").Merge()
$ws.Range("A5:A7").Merge()
$ws.Range("B8:C8").Merge()
$ws.Range("B9:C9").Merge()
$ws.Range("B10:C10").Merge()
$ws.Range("A8:A10").Merge()
$ws.Range("B11:C11").Merge()
$ws.Range("A12:A13").Merge()
$ws.Range("A14:A15").Merge()
$ws.Range("A16:A18").Merge()
$ws.Range("A19:A21").Merge()
$ws.Range("B22:C22").Merge()
$ws.Range("B23:C23").Merge()
$ws.Range("B24:C24").Merge()
$ws.Range("A22:A24").Merge()
$ws.Range("A25:A28").Merge()
